$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new direction-suffixed "N" (near-door) clue strings by setting
# these specific cells to new values. Excel will add them to the shared
# string table as new unique strings since they don't currently exist.
$ws.Range("A13").Value = "KN"
$ws.Range("A3").Value  = "GN"
$ws.Range("G3").Value  = "RN"
$ws.Range("N3").Value  = "YN"
$ws.Range("U5").Value  = "LN"
$ws.Range("U18").Value = "ON"
$ws.Range("M22").Value = "BN"
$ws.Range("G23").Value = "PN"
$ws.Range("A22").Value = "SN"

# Move the active selection to A22 (as saved in the workbook view).
$ws.Range("A22").Select()
